$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range('B14').Value = 6772177
$ws.Range('E14').Value = 'Aguilas Doradas'
$ws.Range('F14').Value = 'Alianza Petrolera'
$ws.Range('G14').Value = 1
$ws.Range('H14').Value = 1
$ws.Range('I14').Value = 'D'
$ws.Range('J14').Value = 2.15
$ws.Range('K14').Value = 3.3
$ws.Range('L14').Value = 3.5
$ws.Range('M14').Value = 2.2
$ws.Range('N14').Value = 3.5
$ws.Range('O14').Value = 3.2
$ws.Range('P14').Value = -0.25
$ws.Range('Q14').Value = 1.9
$ws.Range('R14').Value = 1.9
$ws.Range('S14').Value = 2.75
$ws.Range('T14').Value = 1.95
$ws.Range('U14').Value = 1.85
$ws.Range('V14').Value = -1
$ws.Range('W14').Value = 2.5
$ws.Range('Y14').Value = -0.5
$ws.Range('Z14').Value = 0.45
$ws.Range('AA14').Value = -1
$ws.Range('AB14').Value = 0.8500000000000001
# Row 15
$ws.Range('B15').Value = 6772175
$ws.Range('E15').Value = 'Atletico Nacional Medellin'
$ws.Range('F15').Value = 'Deportivo Pasto'
$ws.Range('G15').Value = 3
$ws.Range('H15').Value = 2
$ws.Range('I15').Value = 'H'
$ws.Range('J15').Value = 1.666
$ws.Range('K15').Value = 3.75
$ws.Range('L15').Value = 4.5
$ws.Range('M15').Value = 1.8
$ws.Range('N15').Value = 3.6
$ws.Range('O15').Value = 5
$ws.Range('P15').Value = -0.75
$ws.Range('Q15').Value = 2
$ws.Range('R15').Value = 1.85
$ws.Range('S15').Value = 2.25
$ws.Range('T15').Value = 1.85
$ws.Range('U15').Value = 2
$ws.Range('V15').Value = 0.8
$ws.Range('W15').Value = -1
$ws.Range('Y15').Value = 0.5
$ws.Range('Z15').Value = -0.5
$ws.Range('AA15').Value = 0.8500000000000001
$ws.Range('AB15').Value = -1
# Row 208
$ws.Range('B208').Value = 7404218
$ws.Range('E208').Value = 'Junior'
$ws.Range('F208').Value = 'Atletico Huila'
$ws.Range('G208').Value = 2
$ws.Range('H208').Value = 0
$ws.Range('I208').Value = 'H'
$ws.Range('J208').Value = 1.363
$ws.Range('K208').Value = 4.5
$ws.Range('L208').Value = 7
$ws.Range('M208').Value = 1.3
$ws.Range('N208').Value = 5
$ws.Range('O208').Value = 12
$ws.Range('P208').Value = -1.5
$ws.Range('Q208').Value = 1.9
$ws.Range('R208').Value = 1.95
$ws.Range('S208').Value = 2.75
$ws.Range('V208').Value = 0.3
$ws.Range('W208').Value = -1
$ws.Range('Y208').Value = 0.8999999999999999
$ws.Range('Z208').Value = -1
# Row 209
$ws.Range('B209').Value = 7404217
$ws.Range('E209').Value = 'Alianza Petrolera'
$ws.Range('F209').Value = 'Deportivo Pereira'
$ws.Range('G209').Value = 2
$ws.Range('I209').Value = 'H'
$ws.Range('J209').Value = 1.95
$ws.Range('K209').Value = 3.2
$ws.Range('L209').Value = 3.75
$ws.Range('M209').Value = 1.95
$ws.Range('N209').Value = 3.2
$ws.Range('O209').Value = 4.75
$ws.Range('P209').Value = -0.5
$ws.Range('Q209').Value = 1.925
$ws.Range('R209').Value = 1.875
$ws.Range('S209').Value = 2
$ws.Range('T209').Value = 1.825
$ws.Range('V209').Value = 0.95
$ws.Range('W209').Value = -1
$ws.Range('Y209').Value = 0.925
$ws.Range('Z209').Value = -1
$ws.Range('AA209').Value = 0.825
$ws.Range('AB209').Value = -1
# Row 211
$ws.Range('B211').Value = 7404214
$ws.Range('E211').Value = 'Boyaca Chico'
$ws.Range('F211').Value = 'Deportivo Cali'
$ws.Range('G211').Value = 1
$ws.Range('I211').Value = 'D'
$ws.Range('J211').Value = 3.2
$ws.Range('K211').Value = 3.1
$ws.Range('L211').Value = 2.2
$ws.Range('M211').Value = 3.6
$ws.Range('N211').Value = 3
$ws.Range('O211').Value = 2.25
$ws.Range('P211').Value = 0.25
$ws.Range('Q211').Value = 1.95
$ws.Range('R211').Value = 1.9
$ws.Range('S211').Value = 2.25
$ws.Range('T211').Value = 1.875
$ws.Range('V211').Value = -1
$ws.Range('W211').Value = 2
$ws.Range('Y211').Value = 0.475
$ws.Range('Z211').Value = -0.5
$ws.Range('AA211').Value = -0.5
$ws.Range('AB211').Value = 0.4875
# Row 212
$ws.Range('B212').Value = 7404212
$ws.Range('E212').Value = 'Envigado FC'
$ws.Range('F212').Value = 'Deportivo Pasto'
$ws.Range('G212').Value = 1
$ws.Range('H212').Value = 1
$ws.Range('I212').Value = 'D'
$ws.Range('J212').Value = 2.6
$ws.Range('K212').Value = 2.875
$ws.Range('L212').Value = 2.8
$ws.Range('M212').Value = 2.8
$ws.Range('N212').Value = 3.2
$ws.Range('O212').Value = 2.625
$ws.Range('P212').Value = 0
$ws.Range('Q212').Value = 1.975
$ws.Range('R212').Value = 1.875
$ws.Range('S212').Value = 2.5
$ws.Range('V212').Value = -1
$ws.Range('W212').Value = 2.2
$ws.Range('Y212').Value = 0
$ws.Range('Z212').Value = 0
# Row 213
$ws.Range('B213').Value = 7404215
$ws.Range('E213').Value = 'America de Cali'
$ws.Range('F213').Value = 'Atletico Bucaramanga'
$ws.Range('G213').Value = 1
$ws.Range('H213').Value = 2
$ws.Range('J213').Value = 1.444
$ws.Range('K213').Value = 4.5
$ws.Range('L213').Value = 6
$ws.Range('M213').Value = 1.363
$ws.Range('N213').Value = 5
$ws.Range('O213').Value = 7.5
$ws.Range('P213').Value = -1.25
$ws.Range('Q213').Value = 1.775
$ws.Range('R213').Value = 2.025
$ws.Range('S213').Value = 3
$ws.Range('T213').Value = 1.925
$ws.Range('U213').Value = 1.875
$ws.Range('X213').Value = 6.5
$ws.Range('Z213').Value = 1.025
$ws.Range('AA213').Value = 0
$ws.Range('AB213').Value = 0
# Row 214
$ws.Range('B214').Value = 7404522
$ws.Range('E214').Value = 'La Equidad'
$ws.Range('F214').Value = 'Millonarios'
$ws.Range('H214').Value = 1
$ws.Range('I214').Value = 'H'
$ws.Range('J214').Value = 2.4
$ws.Range('K214').Value = 3.1
$ws.Range('L214').Value = 2.875
$ws.Range('M214').Value = 2.1
$ws.Range('N214').Value = 3.1
$ws.Range('O214').Value = 3.8
$ws.Range('P214').Value = -0.25
$ws.Range('Q214').Value = 1.75
$ws.Range('R214').Value = 2.05
$ws.Range('S214').Value = 2
$ws.Range('T214').Value = 1.85
$ws.Range('U214').Value = 1.95
$ws.Range('V214').Value = 1.1
$ws.Range('X214').Value = -1
$ws.Range('Y214').Value = 0.75
$ws.Range('Z214').Value = -1
$ws.Range('AA214').Value = 0.8500000000000001
# Row 215
$ws.Range('B215').Value = 7404260
$ws.Range('E215').Value = 'Atletico Nacional Medellin'
$ws.Range('F215').Value = 'Deportes Tolima'
$ws.Range('H215').Value = 3
$ws.Range('I215').Value = 'A'
$ws.Range('J215').Value = 2
$ws.Range('K215').Value = 3.25
$ws.Range('L215').Value = 3.5
$ws.Range('M215').Value = 1.75
$ws.Range('N215').Value = 3.6
$ws.Range('O215').Value = 4.75
$ws.Range('P215').Value = -0.75
$ws.Range('Q215').Value = 2
$ws.Range('R215').Value = 1.8
$ws.Range('S215').Value = 2.5
$ws.Range('T215').Value = 2
$ws.Range('U215').Value = 1.8
$ws.Range('V215').Value = -1
$ws.Range('X215').Value = 3.75
$ws.Range('Y215').Value = -1
$ws.Range('Z215').Value = 0.8
$ws.Range('AA215').Value = 1
# Row 216
$ws.Range('B216').Value = 7404213
$ws.Range('E216').Value = 'Jaguares de Cordoba'
$ws.Range('F216').Value = 'Aguilas Doradas'
$ws.Range('H216').Value = 1
$ws.Range('J216').Value = 3.25
$ws.Range('L216').Value = 2.2
$ws.Range('N216').Value = 3.2
$ws.Range('O216').Value = 2.15
$ws.Range('Q216').Value = 1.975
$ws.Range('R216').Value = 1.825
$ws.Range('S216').Value = 2
$ws.Range('T216').Value = 1.75
$ws.Range('U216').Value = 2.05
$ws.Range('X216').Value = 1.15
$ws.Range('Z216').Value = 0.825
$ws.Range('AA216').Value = -1
$ws.Range('AB216').Value = 1.05
# Row 217
$ws.Range('B217').Value = 7404219
$ws.Range('E217').Value = 'Union Magdalena'
$ws.Range('F217').Value = 'Independiente Medellin'
$ws.Range('G217').Value = 0
$ws.Range('H217').Value = 4
$ws.Range('J217').Value = 3
$ws.Range('K217').Value = 3.1
$ws.Range('L217').Value = 2.3
$ws.Range('M217').Value = 3.6
$ws.Range('N217').Value = 3.4
$ws.Range('O217').Value = 2.1
$ws.Range('P217').Value = 0.25
$ws.Range('Q217').Value = 2.025
$ws.Range('R217').Value = 1.775
$ws.Range('S217').Value = 2.5
$ws.Range('T217').Value = 1.85
$ws.Range('U217').Value = 1.95
$ws.Range('X217').Value = 1.1
$ws.Range('Z217').Value = 0.7749999999999999
$ws.Range('AA217').Value = 0.8500000000000001
$ws.Range('AB217').Value = -1
# Row 240
$ws.Range('B240').Value = 7528603
$ws.Range('E240').Value = 'Junior'
$ws.Range('F240').Value = 'Deportes Tolima'
$ws.Range('G240').Value = 4
$ws.Range('H240').Value = 2
$ws.Range('J240').Value = 1.95
$ws.Range('L240').Value = 4
$ws.Range('M240').Value = 1.909
$ws.Range('N240').Value = 3.75
$ws.Range('O240').Value = 3.8
$ws.Range('P240').Value = -0.5
$ws.Range('Q240').Value = 1.9
$ws.Range('R240').Value = 1.9
$ws.Range('T240').Value = 1.85
$ws.Range('U240').Value = 1.95
$ws.Range('V240').Value = 0.909
$ws.Range('Y240').Value = 0.8999999999999999
$ws.Range('AA240').Value = 0.8500000000000001
# Row 241
$ws.Range('B241').Value = 7528135
$ws.Range('E241').Value = 'Independiente Medellin'
$ws.Range('F241').Value = 'America de Cali'
$ws.Range('G241').Value = 2
$ws.Range('H241').Value = 1
$ws.Range('J241').Value = 2.15
$ws.Range('L241').Value = 3.4
$ws.Range('M241').Value = 2.375
$ws.Range('N241').Value = 3.3
$ws.Range('O241').Value = 3.1
$ws.Range('P241').Value = -0.25
$ws.Range('Q241').Value = 2
$ws.Range('R241').Value = 1.8
$ws.Range('T241').Value = 1.975
$ws.Range('U241').Value = 1.825
$ws.Range('V241').Value = 1.375
$ws.Range('Y241').Value = 1
$ws.Range('AA241').Value = 0.9750000000000001
# Row 430
$ws.Range('M430').Value = 2
$ws.Range('N430').Value = 3.3
$ws.Range('O430').Value = 3.8
$ws.Range('Q430').Value = 2.05
$ws.Range('R430').Value = 1.8
$ws.Range('T430').Value = 1.8
$ws.Range('U430').Value = 2.05
# Row 431
$ws.Range('M431').Value = 1.4
$ws.Range('N431').Value = 4.5
$ws.Range('O431').Value = 8.5
$ws.Range('Q431').Value = 1.9
$ws.Range('R431').Value = 1.95
$ws.Range('T431').Value = 1.975
$ws.Range('U431').Value = 1.875
